$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# New source files being handed off:
#   5685cf38-4c69-4098-a2ff-8993427d9e74.md
#   9d5a481b-0e7e-4d35-83d6-710f7fc0165e.md
# Both reach "Ready for handoff" status with handoff xliffs generated for
# zh-cn and de-de locales.
# ---------------------------------------------------------------------------

$file1Name = "5685cf38-4c69-4098-a2ff-8993427d9e74.md"
$file1Path = "e2e\5685cf38-4c69-4098-a2ff-8993427d9e74.md"
$file1Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a2deea117ce5575e6c705dd36c7155a998422af/e2e/5685cf38-4c69-4098-a2ff-8993427d9e74.md"

$file2Name = "9d5a481b-0e7e-4d35-83d6-710f7fc0165e.md"
$file2Path = "e2e\9d5a481b-0e7e-4d35-83d6-710f7fc0165e.md"
$file2Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a2deea117ce5575e6c705dd36c7155a998422af/e2e/9d5a481b-0e7e-4d35-83d6-710f7fc0165e.md"

$status       = "Ready for handoff"
$overviewDate = "2016-10-14 07:30:34"

$zhXlf1  = "5685cf38-4c69-4098-a2ff-8993427d9e74.cb39d65c2879de697610e859672b22aa6064c84e.zh-cn.xlf"
$zhXlf2  = "9d5a481b-0e7e-4d35-83d6-710f7fc0165e.8cb48b08132462508919a768d8a715353e88301a.zh-cn.xlf"
$zhDate  = "2016-10-14 07:30:24"

$deXlf1  = "5685cf38-4c69-4098-a2ff-8993427d9e74.cb39d65c2879de697610e859672b22aa6064c84e.de-de.xlf"
$deXlf2  = "9d5a481b-0e7e-4d35-83d6-710f7fc0165e.8cb48b08132462508919a768d8a715353e88301a.de-de.xlf"
$deDate  = "2016-10-14 07:30:34"

# NOTE: the source data model stores "True"/"False"/"" as plain text (shared
# strings), matching the existing rows. Writing a bare "False"/"True" value
# through COM auto-coerces to a real Boolean cell, and a bare "" clears the
# cell entirely instead of writing empty text - so a leading apostrophe is
# used here to force plain-text interpretation, exactly like a user typing
# '0001-01-01 style text into a cell in real Excel.
$falseTxt = "'False"
$trueTxt  = "'True"
$emptyTxt = "'"

# ---------------------------------------------------------------------------
# Overview sheet: rows 4 and 5 (columns A-G)
# ---------------------------------------------------------------------------

$overview.Cells.Item(4,1).Value = $file1Name
$overview.Cells.Item(4,2).Value = $file1Path
$overview.Cells.Item(4,3).Value = ".md"
$overview.Cells.Item(4,4).Value = $emptyTxt
$overview.Cells.Item(4,5).Value = $status
$overview.Cells.Item(4,6).Value = $status
$overview.Cells.Item(4,7).Value = $overviewDate

$overview.Cells.Item(5,1).Value = $file2Name
$overview.Cells.Item(5,2).Value = $file2Path
$overview.Cells.Item(5,3).Value = ".md"
$overview.Cells.Item(5,4).Value = $emptyTxt
$overview.Cells.Item(5,5).Value = $status
$overview.Cells.Item(5,6).Value = $status
$overview.Cells.Item(5,7).Value = $overviewDate

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 4 and 5 (columns A-P)
# ---------------------------------------------------------------------------

$zhcn.Cells.Item(4,1).Value  = $file1Name
$zhcn.Cells.Item(4,2).Value  = ".md"
$zhcn.Cells.Item(4,3).Value  = $status
$zhcn.Cells.Item(4,4).Value  = "e2e"
$zhcn.Cells.Item(4,5).Value  = "ht"
$zhcn.Cells.Item(4,6).Value  = $falseTxt
$zhcn.Cells.Item(4,7).Value  = $zhXlf1
$zhcn.Cells.Item(4,8).Value  = $zhDate
$zhcn.Cells.Item(4,9).Value  = $emptyTxt
$zhcn.Cells.Item(4,10).Value = $emptyTxt
$zhcn.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$zhcn.Cells.Item(4,12).Value = $emptyTxt
$zhcn.Cells.Item(4,13).Value = $trueTxt
$zhcn.Cells.Item(4,14).Value = $emptyTxt
$zhcn.Cells.Item(4,15).Value = $falseTxt
$zhcn.Cells.Item(4,16).Value = $emptyTxt

$zhcn.Cells.Item(5,1).Value  = $file2Name
$zhcn.Cells.Item(5,2).Value  = ".md"
$zhcn.Cells.Item(5,3).Value  = $status
$zhcn.Cells.Item(5,4).Value  = "e2e"
$zhcn.Cells.Item(5,5).Value  = "ht"
$zhcn.Cells.Item(5,6).Value  = $falseTxt
$zhcn.Cells.Item(5,7).Value  = $zhXlf2
$zhcn.Cells.Item(5,8).Value  = $zhDate
$zhcn.Cells.Item(5,9).Value  = $emptyTxt
$zhcn.Cells.Item(5,10).Value = $emptyTxt
$zhcn.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$zhcn.Cells.Item(5,12).Value = $emptyTxt
$zhcn.Cells.Item(5,13).Value = $trueTxt
$zhcn.Cells.Item(5,14).Value = $emptyTxt
$zhcn.Cells.Item(5,15).Value = $falseTxt
$zhcn.Cells.Item(5,16).Value = $emptyTxt

# ---------------------------------------------------------------------------
# de-de sheet: rows 4 and 5 (columns A-P)
# ---------------------------------------------------------------------------

$dede.Cells.Item(4,1).Value  = $file1Name
$dede.Cells.Item(4,2).Value  = ".md"
$dede.Cells.Item(4,3).Value  = $status
$dede.Cells.Item(4,4).Value  = "e2e"
$dede.Cells.Item(4,5).Value  = "ht"
$dede.Cells.Item(4,6).Value  = $falseTxt
$dede.Cells.Item(4,7).Value  = $deXlf1
$dede.Cells.Item(4,8).Value  = $deDate
$dede.Cells.Item(4,9).Value  = $emptyTxt
$dede.Cells.Item(4,10).Value = $emptyTxt
$dede.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$dede.Cells.Item(4,12).Value = $emptyTxt
$dede.Cells.Item(4,13).Value = $trueTxt
$dede.Cells.Item(4,14).Value = $emptyTxt
$dede.Cells.Item(4,15).Value = $falseTxt
$dede.Cells.Item(4,16).Value = $emptyTxt

$dede.Cells.Item(5,1).Value  = $file2Name
$dede.Cells.Item(5,2).Value  = ".md"
$dede.Cells.Item(5,3).Value  = $status
$dede.Cells.Item(5,4).Value  = "e2e"
$dede.Cells.Item(5,5).Value  = "ht"
$dede.Cells.Item(5,6).Value  = $falseTxt
$dede.Cells.Item(5,7).Value  = $deXlf2
$dede.Cells.Item(5,8).Value  = $deDate
$dede.Cells.Item(5,9).Value  = $emptyTxt
$dede.Cells.Item(5,10).Value = $emptyTxt
$dede.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$dede.Cells.Item(5,12).Value = $emptyTxt
$dede.Cells.Item(5,13).Value = $trueTxt
$dede.Cells.Item(5,14).Value = $emptyTxt
$dede.Cells.Item(5,15).Value = $falseTxt
$dede.Cells.Item(5,16).Value = $emptyTxt

# ---------------------------------------------------------------------------
# Hyperlinks for the new rows
# ---------------------------------------------------------------------------

$overview.Hyperlinks.Add($overview.Range("B4"), $file1Url, [Type]::Missing, [Type]::Missing, $file1Path) | Out-Null
$overview.Hyperlinks.Add($overview.Range("B5"), $file2Url, [Type]::Missing, [Type]::Missing, $file2Path) | Out-Null

$zhcn.Hyperlinks.Add($zhcn.Range("A4"), $file1Url, [Type]::Missing, [Type]::Missing, $file1Name) | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A5"), $file2Url, [Type]::Missing, [Type]::Missing, $file2Name) | Out-Null

$dede.Hyperlinks.Add($dede.Range("A4"), $file1Url, [Type]::Missing, [Type]::Missing, $file1Name) | Out-Null
$dede.Hyperlinks.Add($dede.Range("A5"), $file2Url, [Type]::Missing, [Type]::Missing, $file2Name) | Out-Null

# ---------------------------------------------------------------------------
# Extend the Excel Tables / AutoFilters to cover the new rows
# ---------------------------------------------------------------------------

$overview.ListObjects.Item(1).Resize($overview.Range("A1:G5"))
$zhcn.ListObjects.Item(1).Resize($zhcn.Range("A1:P5"))
$dede.ListObjects.Item(1).Resize($dede.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Column width adjustments (widened to fit the new status text)
# ---------------------------------------------------------------------------

$overview.Columns.Item(5).ColumnWidth = 17.2159881591797
$overview.Columns.Item(6).ColumnWidth = 17.2159881591797

$zhcn.Columns.Item(3).ColumnWidth = 17.2159881591797
$dede.Columns.Item(3).ColumnWidth = 17.2159881591797
